$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1. Update the "Solo Champ Pool and Play Rate" table (Table 2).
#    Rows 3-7 (1-based) hold Champion / Games / PlayRate.
# ---------------------------------------------------------------
$t = $d.Tables(2)

# Row 3: Aatrox -> Ornn, 21 -> 23, 18% -> 14%
$t.Cell(3,1).Range.Text = "Ornn"
$t.Cell(3,2).Range.Text = "23"
$t.Cell(3,3).Range.Text = "14%"

# Row 4: Ornn -> Aatrox, 18 -> 21, 15% -> 13%
$t.Cell(4,1).Range.Text = "Aatrox"
$t.Cell(4,2).Range.Text = "21"
$t.Cell(4,3).Range.Text = "13%"

# Row 5: Sett -> Gangplank, 13 -> 20, 11% -> 12%
$t.Cell(5,1).Range.Text = "Gangplank"
$t.Cell(5,2).Range.Text = "20"
$t.Cell(5,3).Range.Text = "12%"

# Row 6: Renekton -> Gnar, 11 -> 13, 9% -> 8%
$t.Cell(6,1).Range.Text = "Gnar"
$t.Cell(6,2).Range.Text = "13"
$t.Cell(6,3).Range.Text = "8%"

# Row 7: Gangplank -> Sett, 10 -> 13, 8% unchanged
$t.Cell(7,1).Range.Text = "Sett"
$t.Cell(7,2).Range.Text = "13"

# ---------------------------------------------------------------
# 2. Rewrite the "Solo Champ Pool" discussion paragraph.
#    The original single paragraph (starting "I cut Solo off at
#    top 5...") is split into three paragraphs:
#      (a) new discussion text ending with a _GoBack bookmark
#      (b) an empty paragraph
#      (c) the unchanged bold "This probably explains..." text
# ---------------------------------------------------------------
$p = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs($i)
    if ($cand.Range.Text.StartsWith("I cut Solo off at top")) {
        $p = $cand
        break
    }
}

$pStart = $p.Range.Start
$pEnd = $p.Range.End

# Locate the point where the unchanged bold portion begins.
$splitR = $d.Range($pStart, $pEnd)
$splitR.Find.Execute("This probably explains", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $splitR.Start

# Remove all of the old text (and the old bookmark that lived inside it)
# before that point.
$deleteRange = $d.Range($pStart, $splitPoint)
$deleteRange.Text = ""

# Split off two new empty paragraphs ahead of the (now-shortened)
# "This probably explains..." paragraph.
$breakPos = $d.Range($pStart, $pStart)
$breakPos.InsertParagraphBefore()
$breakPos.InsertParagraphBefore()

# Fill the first of the two new paragraphs with the new discussion text.
$run1 = "Solo’s top 5 most played champions is about what I would expect from watching him over the years. It’s the most popular toplane tank (Ornn), a set of bruisers (Aatrox, Sett, Gnar) and "
$run2 = "a nice blind-pickable carry (Gangplank). This is a pretty st"
$run3 = "andard collection of top lane champions."

$cursor = $d.Range($pStart, $pStart)
$cursor.InsertAfter($run1)
$cursor.Collapse(0)
$cursor.InsertAfter($run2)
$cursor.Collapse(0)
$cursor.InsertAfter($run3)
$cursor.Collapse(0)

$newTextLen = $run1.Length + $run2.Length + $run3.Length
$firstNewParaEnd = $pStart + $newTextLen

$firstPara = $d.Range($pStart, $pStart).Paragraphs(1)
$secondPara = $d.Range($firstNewParaEnd + 1, $firstNewParaEnd + 1).Paragraphs(1)
$firstPara.Range.Font.Bold = 0
$secondPara.Range.Font.Bold = 0

# Re-create the _GoBack bookmark at the end of the new text (it previously
# sat inside the deleted text).  Anchor on the last character, since a
# zero-length range is not accepted reliably.
$bmTarget = $d.Range($firstNewParaEnd - 1, $firstNewParaEnd)
$d.Bookmarks.Add("_GoBack", $bmTarget)

Write-Output "Paragraph restructuring complete"
